$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the same header style/format
# as the neighboring "sum" header in G1 (bold, centered, bordered).
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2 for the single data row.
$ws.Range("H2").Value = 0
